$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "cryptos" price/volume table with the latest scraped values.
# Price cells (column D) whose new value would otherwise be re-interpreted
# by Excel as a plain number (losing formatting such as trailing zeros,
# e.g. "0.100" or "35.20") are written with a leading apostrophe so they
# stay stored as text, matching the original workbook's inlineStr cells.
$ws.Range('D2').Value = '68.468.31'
$ws.Range('E2').Value = '  +1.04%  '
$ws.Range('D3').Value = '3.801.51'
$ws.Range('E3').Value = '  -0.59%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'606.51"
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').Value = "'164.51"
$ws.Range('E6').Value = '  -1.19%  '
$ws.Range('D7').Value = '3.798.66'
$ws.Range('E7').Value = '  -0.59%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').Value = "'0.518"
$ws.Range('E9').Value = '  -0.27%  '
$ws.Range('D10').Value = "'0.159"
$ws.Range('E10').Value = '  -0.85%  '
$ws.Range('D11').Value = "'6.95"
$ws.Range('E11').Value = '  +10.34%  '
$ws.Range('D12').Value = "'0.451"
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').Value = "'0.0000248"
$ws.Range('E13').Value = '  -2.12%  '
$ws.Range('D14').Value = "'35.20"
$ws.Range('E14').Value = '  -2.12%  '
$ws.Range('D15').Value = '4.436.67'
$ws.Range('E15').Value = '  -0.79%  '
$ws.Range('D16').Value = '3.788.47'
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('D17').Value = '68.351.21'
$ws.Range('E17').Value = '  +0.82%  '
$ws.Range('D18').Value = "'18.14"
$ws.Range('E18').Value = '  -1.94%  '
$ws.Range('E19').Value = '  +1.86%  '
$ws.Range('D20').Value = "'7.06"
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('D21').Value = "'462.12"
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('D22').Value = "'9.61"
$ws.Range('E22').Value = '  -3.06%  '
$ws.Range('D23').Value = "'0.698"
$ws.Range('E23').Value = '  -0.64%  '
$ws.Range('D24').Value = "'0.0000150"
$ws.Range('E24').Value = '  +0.56%  '
$ws.Range('D25').Value = "'83.69"
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('D26').Value = "'11.99"
$ws.Range('E26').Value = '  -1.35%  '
$ws.Range('B27').Value = 'Fetch.AI'
$ws.Range('C27').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D27').Value = "'2.11"
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').Value = "'10.05"
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').Value = "'1.00"
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').Value = '3.943.64'
$ws.Range('E30').Value = '  -0.83%  '
$ws.Range('D31').Value = "'2.62"
$ws.Range('E31').Value = '  -6.09%  '
$ws.Range('D32').Value = "'7.26"
$ws.Range('E32').Value = '  -1.48%  '
$ws.Range('D33').Value = "'2.20"
$ws.Range('E33').Value = '  -1.35%  '
$ws.Range('D34').Value = "'29.14"
$ws.Range('E34').Value = '  -1.70%  '
$ws.Range('D35').Value = "'0.997"
$ws.Range('E35').Value = '  -0.36%  '
$ws.Range('D36').Value = "'8.99"
$ws.Range('E36').Value = '  -1.41%  '
$ws.Range('D37').Value = "'0.100"
$ws.Range('E37').Value = '  +0.30%  '
$ws.Range('D38').Value = "'0.150"
$ws.Range('E38').Value = '  +8.44%  '
$ws.Range('D39').Value = "'5.89"
$ws.Range('E39').Value = '  +1.19%  '
$ws.Range('D40').Value = "'3.26"
$ws.Range('E40').Value = '  -1.28%  '
$ws.Range('D41').Value = "'0.982"
$ws.Range('E41').Value = '  -1.47%  '
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('B44').Value = 'ONDO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D44').Value = "'1.43"
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('B45').Value = 'Arweave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D45').Value = "'43.55"
$ws.Range('E45').Value = '  -2.95%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = "'47.05"
$ws.Range('E46').Value = '  -1.32%  '
$ws.Range('B47').Value = 'TheGraph'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D47').Value = "'0.296"
$ws.Range('E47').Value = '  -1.49%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = "'152.90"
$ws.Range('E48').Value = '  +1.24%  '
$ws.Range('D49').Value = "'8.38"
$ws.Range('E49').Value = '  +0.18%  '
$ws.Range('D50').Value = "'1.86"
$ws.Range('E50').Value = '  +0.22%  '
$ws.Range('D51').Value = "'26.30"
$ws.Range('E51').Value = '  -8.84%  '
